# This edit corresponds to an XML re-serialization change only (attribute
# ordering of namespace declarations and element attributes, caused by an
# Apache POI library upgrade). There is no semantic change to the document's
# content, formatting, or structure, so no Word object-model operations are
# required here.
$d = $word.ActiveDocument
